$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.2272265
$ws.Range("N2").Value = 0.454453
$ws.Range("O2").Value = 0.08704083604617911
$ws.Range("P2").Value = 0.08229687998280369
$ws.Range("Q2").Value = 0.130550258857
$ws.Range("R2").Value = 0.783301553142
$ws.Range("S2").Value = 0.08704083604617911
$ws.Range("T2").Value = 0.08229687998280369

# Row 3
$ws.Range("O3").Value = 0.1052353694185077
$ws.Range("P3").Value = 0.149249644656207
$ws.Range("S3").Value = 0.1052353694185077
$ws.Range("T3").Value = 0.149249644656207

# Row 4
$ws.Range("M4").Value = 2.082377
$ws.Range("N4").Value = 4.164754
$ws.Range("O4").Value = 0.7976703203338269
$ws.Range("P4").Value = 0.7541951755096822
$ws.Range("Q4").Value = 1.196404716826
$ws.Range("R4").Value = 7.178428300956001
$ws.Range("S4").Value = 0.7976703203338269
$ws.Range("T4").Value = 0.7541951755096822

# Row 5
$ws.Range("M5").Value = 0.02610733333333333
$ws.Range("N5").Value = 0.078322
$ws.Range("O5").Value = 0.0100006122537187
$ws.Range("P5").Value = 0.01418332860386696
$ws.Range("Q5").Value = 0.01499965507866667
$ws.Range("R5").Value = 0.134996895708
$ws.Range("S5").Value = 0.0100006122537187
$ws.Range("T5").Value = 0.01418332860386696

# Row 6
$ws.Range("M6").Value = 0.000138
$ws.Range("N6").Value = 0.000414
$ws.Range("O6").Value = 0.00005286194776741585
$ws.Range("P6").Value = 0.0000749712474400669
$ws.Range("Q6").Value = 0.000079286244
$ws.Range("R6").Value = 0.0007135761959999999
$ws.Range("S6").Value = 0.00005286194776741585
$ws.Range("T6").Value = 0.0000749712474400669
